$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: tickers for the new companies (rows 96-140)
$tickers = @(
    ,@(96, 'TM')
    ,@(97, 'SAP')
    ,@(98, 'BABA')
    ,@(99, 'NSRGY')
    ,@(100, 'SNY')
    ,@(101, 'ORLY')
    ,@(102, 'UPS')
    ,@(103, 'CVS')
    ,@(104, 'ADP')
    ,@(105, 'C')
    ,@(106, 'MELI')
    ,@(107, 'F')
    ,@(108, 'INTU')
    ,@(109, 'DELL')
    ,@(110, 'PYPL')
    ,@(111, 'ABNB')
    ,@(112, 'BK')
    ,@(113, 'EXC')
    ,@(114, 'MU')
    ,@(115, 'EQIX')
    ,@(116, 'ILMN')
    ,@(117, 'PLTR')
    ,@(118, 'WDAY')
    ,@(119, 'TWTR')
    ,@(120, 'SNOW')
    ,@(121, 'ZM')
    ,@(122, 'DOCU')
    ,@(123, 'ROST')
    ,@(124, 'VZ')
    ,@(125, 'DLR')
    ,@(126, 'HCA')
    ,@(127, 'FISV')
    ,@(128, 'MAR')
    ,@(129, 'ECL')
    ,@(130, 'LHX')
    ,@(131, 'CTSH')
    ,@(132, 'CCI')
    ,@(133, 'FTNT')
    ,@(134, 'NOC')
    ,@(135, 'WBA')
    ,@(136, 'DHR')
    ,@(137, 'TROW')
    ,@(138, 'ALGN')
    ,@(139, 'GIS')
    ,@(140, 'NEM')
)
foreach ($t in $tickers) {
    $ws.Cells.Item($t[0], 1).Value = $t[1]
}

# Column B: company names for the new companies (rows 96-142)
$companyNames = @(
    ,@(96, 'Toyota Motor Corporation  ')
    ,@(97, 'SAP SE  ')
    ,@(98, 'Alibaba Group  ')
    ,@(99, 'Nestlé  ')
    ,@(100, 'Sanofi  ')
    ,@(101, 'O''Reilly Automotive  ')
    ,@(102, 'United Parcel Service  ')
    ,@(103, 'CVS Health  ')
    ,@(104, 'Automatic Data Processing  ')
    ,@(105, 'Citigroup  ')
    ,@(106, 'MercadoLibre  ')
    ,@(107, 'Ford Motor Company  ')
    ,@(108, 'Intuit  ')
    ,@(109, 'Dell Technologies  ')
    ,@(110, 'PayPal Holdings  ')
    ,@(111, 'Airbnb  ')
    ,@(112, 'Bank of New York Mellon  ')
    ,@(113, 'Exelon Corporation  ')
    ,@(114, 'Micron Technology  ')
    ,@(115, 'Equinix  ')
    ,@(116, 'Intuitive Surgical  ')
    ,@(117, 'Illumina  ')
    ,@(118, 'Palantir Technologies  ')
    ,@(119, 'Workday  ')
    ,@(120, 'Twitter (X)  ')
    ,@(121, 'Snowflake Inc.  ')
    ,@(122, 'Zoom Video Communications  ')
    ,@(123, 'DocuSign  ')
    ,@(124, 'Ross Stores  ')
    ,@(125, 'Verizon Communications  ')
    ,@(126, 'Digital Realty Trust  ')
    ,@(127, 'HCA Healthcare  ')
    ,@(128, 'Fiserv  ')
    ,@(129, 'Marriott International  ')
    ,@(130, 'Ecolab  ')
    ,@(131, 'Sherwin-Williams  ')
    ,@(132, 'L3Harris Technologies  ')
    ,@(133, 'Cognizant Technology Solutions  ')
    ,@(134, 'Crown Castle International  ')
    ,@(135, 'Fortinet  ')
    ,@(136, 'Northrop Grumman  ')
    ,@(137, 'Walgreens Boots Alliance  ')
    ,@(138, 'Danaher Corporation  ')
    ,@(139, 'T. Rowe Price Group  ')
    ,@(140, 'Align Technology  ')
    ,@(141, 'General Mills  ')
    ,@(142, 'Newmont Corporation ')
)
foreach ($c in $companyNames) {
    $ws.Cells.Item($c[0], 2).Value = $c[1]
}

# Update the saved view state (scroll position + active selection)
$ws.Range("F30").Select()
try {
    $excel.ActiveWindow.ScrollRow = 25
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

